$wb = $excel.ActiveWorkbook

# --- ALERTS: append rows 2-3 ---
$ws = $wb.Worksheets.Item("ALERTS")
$ws.Range("A2:A3").NumberFormat = "@"

$ws.Cells.Item(2,1).Value = "2026-02-04"
$ws.Cells.Item(2,2).Value = "14:12:08"
$ws.Cells.Item(2,3).Value = "14:00"
$ws.Cells.Item(2,4).Value = "Bathroom"
$ws.Cells.Item(2,5).Value = "MINIMAL"
$ws.Cells.Item(2,6).Value = "MINIMAL ALERT: Bathroom occupied, no motion > 20s."

$ws.Cells.Item(3,1).Value = "2026-02-04"
$ws.Cells.Item(3,2).Value = "14:12:20"
$ws.Cells.Item(3,3).Value = "14:00"
$ws.Cells.Item(3,4).Value = "Bathroom"
$ws.Cells.Item(3,5).Value = "MODERATE"
$ws.Cells.Item(3,6).Value = "MODERATE ALERT: Bathroom occupied, no motion > 40s."

# --- PIR: append rows 126-139 ---
$ws = $wb.Worksheets.Item("PIR")
$ws.Range("A126:A139").NumberFormat = "@"

$ws.Cells.Item(126,1).Value = "2026-02-04"
$ws.Cells.Item(126,2).Value = "14:11:28"
$ws.Cells.Item(126,3).Value = "14:00"
$ws.Cells.Item(126,4).Value = "Bathroom"
$ws.Cells.Item(126,5).Value = "No Motion"
$ws.Cells.Item(126,6).Value = "Inactive"

$ws.Cells.Item(127,1).Value = "2026-02-04"
$ws.Cells.Item(127,2).Value = "14:11:29"
$ws.Cells.Item(127,3).Value = "14:00"
$ws.Cells.Item(127,4).Value = "Bathroom"
$ws.Cells.Item(127,5).Value = "No Motion"
$ws.Cells.Item(127,6).Value = "Inactive"

$ws.Cells.Item(128,1).Value = "2026-02-04"
$ws.Cells.Item(128,2).Value = "14:11:31"
$ws.Cells.Item(128,3).Value = "14:00"
$ws.Cells.Item(128,4).Value = "Bathroom"
$ws.Cells.Item(128,5).Value = "No Motion"
$ws.Cells.Item(128,6).Value = "Inactive"

$ws.Cells.Item(129,1).Value = "2026-02-04"
$ws.Cells.Item(129,2).Value = "14:11:36"
$ws.Cells.Item(129,3).Value = "14:00"
$ws.Cells.Item(129,4).Value = "Bathroom"
$ws.Cells.Item(129,5).Value = "No Motion"
$ws.Cells.Item(129,6).Value = "Inactive"

$ws.Cells.Item(130,1).Value = "2026-02-04"
$ws.Cells.Item(130,2).Value = "14:11:41"
$ws.Cells.Item(130,3).Value = "14:00"
$ws.Cells.Item(130,4).Value = "Bathroom"
$ws.Cells.Item(130,5).Value = "No Motion"
$ws.Cells.Item(130,6).Value = "Inactive"

$ws.Cells.Item(131,1).Value = "2026-02-04"
$ws.Cells.Item(131,2).Value = "14:11:46"
$ws.Cells.Item(131,3).Value = "14:00"
$ws.Cells.Item(131,4).Value = "Bathroom"
$ws.Cells.Item(131,5).Value = "No Motion"
$ws.Cells.Item(131,6).Value = "Inactive"

$ws.Cells.Item(132,1).Value = "2026-02-04"
$ws.Cells.Item(132,2).Value = "14:11:51"
$ws.Cells.Item(132,3).Value = "14:00"
$ws.Cells.Item(132,4).Value = "Bathroom"
$ws.Cells.Item(132,5).Value = "No Motion"
$ws.Cells.Item(132,6).Value = "Inactive"

$ws.Cells.Item(133,1).Value = "2026-02-04"
$ws.Cells.Item(133,2).Value = "14:11:56"
$ws.Cells.Item(133,3).Value = "14:00"
$ws.Cells.Item(133,4).Value = "Bathroom"
$ws.Cells.Item(133,5).Value = "No Motion"
$ws.Cells.Item(133,6).Value = "Inactive"

$ws.Cells.Item(134,1).Value = "2026-02-04"
$ws.Cells.Item(134,2).Value = "14:12:01"
$ws.Cells.Item(134,3).Value = "14:00"
$ws.Cells.Item(134,4).Value = "Bathroom"
$ws.Cells.Item(134,5).Value = "No Motion"
$ws.Cells.Item(134,6).Value = "Inactive"

$ws.Cells.Item(135,1).Value = "2026-02-04"
$ws.Cells.Item(135,2).Value = "14:12:06"
$ws.Cells.Item(135,3).Value = "14:00"
$ws.Cells.Item(135,4).Value = "Bathroom"
$ws.Cells.Item(135,5).Value = "No Motion"
$ws.Cells.Item(135,6).Value = "Inactive"

$ws.Cells.Item(136,1).Value = "2026-02-04"
$ws.Cells.Item(136,2).Value = "14:12:11"
$ws.Cells.Item(136,3).Value = "14:00"
$ws.Cells.Item(136,4).Value = "Bathroom"
$ws.Cells.Item(136,5).Value = "No Motion"
$ws.Cells.Item(136,6).Value = "Inactive"

$ws.Cells.Item(137,1).Value = "2026-02-04"
$ws.Cells.Item(137,2).Value = "14:12:16"
$ws.Cells.Item(137,3).Value = "14:00"
$ws.Cells.Item(137,4).Value = "Bathroom"
$ws.Cells.Item(137,5).Value = "No Motion"
$ws.Cells.Item(137,6).Value = "Inactive"

$ws.Cells.Item(138,1).Value = "2026-02-04"
$ws.Cells.Item(138,2).Value = "14:12:22"
$ws.Cells.Item(138,3).Value = "14:00"
$ws.Cells.Item(138,4).Value = "Bathroom"
$ws.Cells.Item(138,5).Value = "No Motion"
$ws.Cells.Item(138,6).Value = "Inactive"

$ws.Cells.Item(139,1).Value = "2026-02-04"
$ws.Cells.Item(139,2).Value = "14:12:27"
$ws.Cells.Item(139,3).Value = "14:00"
$ws.Cells.Item(139,4).Value = "Bathroom"
$ws.Cells.Item(139,5).Value = "No Motion"
$ws.Cells.Item(139,6).Value = "Inactive"

# --- Humidity: append rows 101-109 ---
$ws = $wb.Worksheets.Item("Humidity")
$ws.Range("A101:A109").NumberFormat = "@"
$ws.Range("E101:E109").NumberFormat = "@"

$ws.Cells.Item(101,1).Value = "2026-02-04"
$ws.Cells.Item(101,2).Value = "14:11:28"
$ws.Cells.Item(101,3).Value = "14:00"
$ws.Cells.Item(101,4).Value = "Bathroom"
$ws.Cells.Item(101,5).Value = "76.7%"
$ws.Cells.Item(101,6).Value = "Active"

$ws.Cells.Item(102,1).Value = "2026-02-04"
$ws.Cells.Item(102,2).Value = "14:11:35"
$ws.Cells.Item(102,3).Value = "14:00"
$ws.Cells.Item(102,4).Value = "Bathroom"
$ws.Cells.Item(102,5).Value = "76.6%"
$ws.Cells.Item(102,6).Value = "Active"

$ws.Cells.Item(103,1).Value = "2026-02-04"
$ws.Cells.Item(103,2).Value = "14:11:45"
$ws.Cells.Item(103,3).Value = "14:00"
$ws.Cells.Item(103,4).Value = "Bathroom"
$ws.Cells.Item(103,5).Value = "76.5%"
$ws.Cells.Item(103,6).Value = "Active"

$ws.Cells.Item(104,1).Value = "2026-02-04"
$ws.Cells.Item(104,2).Value = "14:11:55"
$ws.Cells.Item(104,3).Value = "14:00"
$ws.Cells.Item(104,4).Value = "Bathroom"
$ws.Cells.Item(104,5).Value = "77.6%"
$ws.Cells.Item(104,6).Value = "Active"

$ws.Cells.Item(105,1).Value = "2026-02-04"
$ws.Cells.Item(105,2).Value = "14:12:00"
$ws.Cells.Item(105,3).Value = "14:00"
$ws.Cells.Item(105,4).Value = "Bathroom"
$ws.Cells.Item(105,5).Value = "77.4%"
$ws.Cells.Item(105,6).Value = "Active"

$ws.Cells.Item(106,1).Value = "2026-02-04"
$ws.Cells.Item(106,2).Value = "14:12:05"
$ws.Cells.Item(106,3).Value = "14:00"
$ws.Cells.Item(106,4).Value = "Bathroom"
$ws.Cells.Item(106,5).Value = "77.4%"
$ws.Cells.Item(106,6).Value = "Active"

$ws.Cells.Item(107,1).Value = "2026-02-04"
$ws.Cells.Item(107,2).Value = "14:12:10"
$ws.Cells.Item(107,3).Value = "14:00"
$ws.Cells.Item(107,4).Value = "Bathroom"
$ws.Cells.Item(107,5).Value = "77.3%"
$ws.Cells.Item(107,6).Value = "Active"

$ws.Cells.Item(108,1).Value = "2026-02-04"
$ws.Cells.Item(108,2).Value = "14:12:21"
$ws.Cells.Item(108,3).Value = "14:00"
$ws.Cells.Item(108,4).Value = "Bathroom"
$ws.Cells.Item(108,5).Value = "77.4%"
$ws.Cells.Item(108,6).Value = "Active"

$ws.Cells.Item(109,1).Value = "2026-02-04"
$ws.Cells.Item(109,2).Value = "14:12:26"
$ws.Cells.Item(109,3).Value = "14:00"
$ws.Cells.Item(109,4).Value = "Bathroom"
$ws.Cells.Item(109,5).Value = "76.3%"
$ws.Cells.Item(109,6).Value = "Active"

# --- Temperature: append rows 101-109 ---
$ws = $wb.Worksheets.Item("Temperature")
$ws.Range("A101:A109").NumberFormat = "@"

$ws.Cells.Item(101,1).Value = "2026-02-04"
$ws.Cells.Item(101,2).Value = "14:11:29"
$ws.Cells.Item(101,3).Value = "14:00"
$ws.Cells.Item(101,4).Value = "Bathroom"
$ws.Cells.Item(101,5).Value = "24.9C"
$ws.Cells.Item(101,6).Value = "Active"

$ws.Cells.Item(102,1).Value = "2026-02-04"
$ws.Cells.Item(102,2).Value = "14:11:36"
$ws.Cells.Item(102,3).Value = "14:00"
$ws.Cells.Item(102,4).Value = "Bathroom"
$ws.Cells.Item(102,5).Value = "24.9C"
$ws.Cells.Item(102,6).Value = "Active"

$ws.Cells.Item(103,1).Value = "2026-02-04"
$ws.Cells.Item(103,2).Value = "14:11:46"
$ws.Cells.Item(103,3).Value = "14:00"
$ws.Cells.Item(103,4).Value = "Bathroom"
$ws.Cells.Item(103,5).Value = "24.9C"
$ws.Cells.Item(103,6).Value = "Active"

$ws.Cells.Item(104,1).Value = "2026-02-04"
$ws.Cells.Item(104,2).Value = "14:11:56"
$ws.Cells.Item(104,3).Value = "14:00"
$ws.Cells.Item(104,4).Value = "Bathroom"
$ws.Cells.Item(104,5).Value = "24.9C"
$ws.Cells.Item(104,6).Value = "Active"

$ws.Cells.Item(105,1).Value = "2026-02-04"
$ws.Cells.Item(105,2).Value = "14:12:01"
$ws.Cells.Item(105,3).Value = "14:00"
$ws.Cells.Item(105,4).Value = "Bathroom"
$ws.Cells.Item(105,5).Value = "24.9C"
$ws.Cells.Item(105,6).Value = "Active"

$ws.Cells.Item(106,1).Value = "2026-02-04"
$ws.Cells.Item(106,2).Value = "14:12:06"
$ws.Cells.Item(106,3).Value = "14:00"
$ws.Cells.Item(106,4).Value = "Bathroom"
$ws.Cells.Item(106,5).Value = "24.9C"
$ws.Cells.Item(106,6).Value = "Active"

$ws.Cells.Item(107,1).Value = "2026-02-04"
$ws.Cells.Item(107,2).Value = "14:12:11"
$ws.Cells.Item(107,3).Value = "14:00"
$ws.Cells.Item(107,4).Value = "Bathroom"
$ws.Cells.Item(107,5).Value = "24.9C"
$ws.Cells.Item(107,6).Value = "Active"

$ws.Cells.Item(108,1).Value = "2026-02-04"
$ws.Cells.Item(108,2).Value = "14:12:21"
$ws.Cells.Item(108,3).Value = "14:00"
$ws.Cells.Item(108,4).Value = "Bathroom"
$ws.Cells.Item(108,5).Value = "25.0C"
$ws.Cells.Item(108,6).Value = "Active"

$ws.Cells.Item(109,1).Value = "2026-02-04"
$ws.Cells.Item(109,2).Value = "14:12:26"
$ws.Cells.Item(109,3).Value = "14:00"
$ws.Cells.Item(109,4).Value = "Bathroom"
$ws.Cells.Item(109,5).Value = "24.9C"
$ws.Cells.Item(109,6).Value = "Active"

# --- Proximity: append rows 5-6 ---
$ws = $wb.Worksheets.Item("Proximity")
$ws.Range("A5:A6").NumberFormat = "@"

$ws.Cells.Item(5,1).Value = "2026-02-04"
$ws.Cells.Item(5,2).Value = "14:11:29"
$ws.Cells.Item(5,3).Value = "14:00"
$ws.Cells.Item(5,4).Value = "Bathroom Door"
$ws.Cells.Item(5,5).Value = "EXIT"
$ws.Cells.Item(5,6).Value = "User EXITED Bathroom"

$ws.Cells.Item(6,1).Value = "2026-02-04"
$ws.Cells.Item(6,2).Value = "14:11:37"
$ws.Cells.Item(6,3).Value = "14:00"
$ws.Cells.Item(6,4).Value = "Bathroom Door"
$ws.Cells.Item(6,5).Value = "ENTER"
$ws.Cells.Item(6,6).Value = "User ENTERED Bathroom"

